# chore: update Sheets via scheduled runner
# Refreshes market-price-derived columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# for a handful of leve rows across the ALC/ARM/BSM/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1595.4
$ws.Range("J17").Value = 1604.2565
$ws.Range("L17").Value = 4812.7695
$ws.Range("N17").Value = -5148.7695
$ws.Range("H98").Value = 11364624
$ws.Range("I98").Value = 15625608
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 15625608
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -15624110
$ws.Range("N98").Value = -4996
$ws.Range("H116").Value = 15848.667
$ws.Range("I116").Value = 6369.25
$ws.Range("J116").Value = 19295.727
$ws.Range("K116").Value = 6369.25
$ws.Range("L116").Value = 19295.727
$ws.Range("M116").Value = -2927.25
$ws.Range("N116").Value = -26179.727
$ws.Range("H122").Value = 11364624
$ws.Range("I122").Value = 15625608
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 46876824
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -46874374
$ws.Range("N122").Value = -10900
$ws.Range("H137").Value = 2372.6667
$ws.Range("I137").Value = 1598.1666
$ws.Range("J137").Value = 3405.3333
$ws.Range("K137").Value = 4794.4998
$ws.Range("L137").Value = 10215.9999
$ws.Range("M137").Value = -2244.4998
$ws.Range("N137").Value = -15315.9999
$ws.Range("H138").Value = 2671.907
$ws.Range("I138").Value = 2842.2778
$ws.Range("J138").Value = 2549.24
$ws.Range("K138").Value = 8526.8334
$ws.Range("L138").Value = 7647.719999999999
$ws.Range("M138").Value = -3386.8334
$ws.Range("N138").Value = -17927.72

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4253.8667
$ws.Range("I122").Value = 4253.8667
$ws.Range("K122").Value = 12761.6001
$ws.Range("M122").Value = -10311.6001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2474.2856
$ws.Range("I86").Value = 1696.6875
$ws.Range("J86").Value = 4962.6
$ws.Range("K86").Value = 1696.6875
$ws.Range("L86").Value = 4962.6
$ws.Range("M86").Value = -573.6875
$ws.Range("N86").Value = -7208.6
$ws.Range("H87").Value = 99949.5
$ws.Range("J87").Value = 99949.5
$ws.Range("L87").Value = 99949.5
$ws.Range("N87").Value = -102445.5
$ws.Range("H89").Value = 2474.2856
$ws.Range("I89").Value = 1696.6875
$ws.Range("J89").Value = 4962.6
$ws.Range("K89").Value = 8483.4375
$ws.Range("L89").Value = 24813
$ws.Range("M89").Value = -2867.4375
$ws.Range("N89").Value = -36045
$ws.Range("H90").Value = 99949.5
$ws.Range("J90").Value = 99949.5
$ws.Range("L90").Value = 299848.5
$ws.Range("N90").Value = -312328.5
$ws.Range("H105").Value = 616295.5600000001
$ws.Range("I105").Value = 859194.25
$ws.Range("J105").Value = 9048.833000000001
$ws.Range("K105").Value = 859194.25
$ws.Range("L105").Value = 9048.833000000001
$ws.Range("M105").Value = -857447.25
$ws.Range("N105").Value = -12542.833

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 200007920
$ws.Range("J42").Value = 12404.333
$ws.Range("L42").Value = 37212.999
$ws.Range("N42").Value = -38280.999
$ws.Range("H113").Value = 1050.5714
$ws.Range("J113").Value = 986.125
$ws.Range("L113").Value = 2958.375
$ws.Range("N113").Value = -7298.375
$ws.Range("H137").Value = 17212.375
$ws.Range("J137").Value = 22999.8
$ws.Range("L137").Value = 68999.39999999999
$ws.Range("N137").Value = -79199.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 28748.75
$ws.Range("I5").Value = 32497.5
$ws.Range("J5").Value = 25000
$ws.Range("K5").Value = 32497.5
$ws.Range("L5").Value = 25000
$ws.Range("M5").Value = -32385.5
$ws.Range("N5").Value = -25224
$ws.Range("H11").Value = 1000
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H70").Value = 11180.52
$ws.Range("I70").Value = 11153.154
$ws.Range("J70").Value = 11189.641
$ws.Range("K70").Value = 11153.154
$ws.Range("L70").Value = 11189.641
$ws.Range("M70").Value = -10883.154
$ws.Range("N70").Value = -11729.641
$ws.Range("H73").Value = 11180.52
$ws.Range("I73").Value = 11153.154
$ws.Range("J73").Value = 11189.641
$ws.Range("K73").Value = 11153.154
$ws.Range("L73").Value = 11189.641
$ws.Range("M73").Value = -10217.154
$ws.Range("N73").Value = -13061.641
$ws.Range("H80").Value = 3184.375
$ws.Range("I80").Value = 2996.5715
$ws.Range("K80").Value = 2996.5715
$ws.Range("M80").Value = -1998.5715
$ws.Range("H83").Value = 3184.375
$ws.Range("I83").Value = 2996.5715
$ws.Range("K83").Value = 14982.8575
$ws.Range("M83").Value = -9990.8575
$ws.Range("H102").Value = 2021.5454
$ws.Range("I102").Value = 2101.2
$ws.Range("J102").Value = 1225
$ws.Range("K102").Value = 2101.2
$ws.Range("L102").Value = 1225
$ws.Range("M102").Value = -479.1999999999998
$ws.Range("N102").Value = -4469

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("H12").Value = 2609.5
$ws.Range("I12").Value = 219
$ws.Range("J12").Value = 5000
$ws.Range("K12").Value = 219
$ws.Range("L12").Value = 5000
$ws.Range("M12").Value = -49
$ws.Range("N12").Value = -5340
$ws.Range("H68").Value = 2454268.5
$ws.Range("I68").Value = 4168076.5
$ws.Range("J68").Value = 5971.143
$ws.Range("K68").Value = 4168076.5
$ws.Range("L68").Value = 5971.143
$ws.Range("M68").Value = -4167327.5
$ws.Range("N68").Value = -7469.143
$ws.Range("H71").Value = 2454268.5
$ws.Range("I71").Value = 4168076.5
$ws.Range("J71").Value = 5971.143
$ws.Range("K71").Value = 20840382.5
$ws.Range("L71").Value = 29855.715
$ws.Range("M71").Value = -20836638.5
$ws.Range("N71").Value = -37343.715
$ws.Range("H136").Value = 5182.385
$ws.Range("I136").Value = 1859.25
$ws.Range("K136").Value = 5577.75
$ws.Range("M136").Value = -3027.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 19003.666
$ws.Range("J10").Value = 11005.5
$ws.Range("L10").Value = 11005.5
$ws.Range("N10").Value = -11343.5
$ws.Range("H13").Value = 4332.3335
$ws.Range("J13").Value = 4998.5
$ws.Range("L13").Value = 4998.5
$ws.Range("N13").Value = -5278.5
$ws.Range("H62").Value = 11562.375
$ws.Range("I62").Value = 4900
$ws.Range("J62").Value = 13783.167
$ws.Range("K62").Value = 4900
$ws.Range("L62").Value = 13783.167
$ws.Range("M62").Value = -4276
$ws.Range("N62").Value = -15031.167
$ws.Range("H65").Value = 11562.375
$ws.Range("I65").Value = 4900
$ws.Range("J65").Value = 13783.167
$ws.Range("K65").Value = 24500
$ws.Range("L65").Value = 68915.83499999999
$ws.Range("M65").Value = -21380
$ws.Range("N65").Value = -75155.83499999999
$ws.Range("H122").Value = 2241.75
$ws.Range("I122").Value = 1802.25
$ws.Range("J122").Value = 3999.75
$ws.Range("K122").Value = 5406.75
$ws.Range("L122").Value = 11999.25
$ws.Range("M122").Value = -2956.75
$ws.Range("N122").Value = -16899.25
